# Auto-generated edit script: bulk update of market-data columns (H-N)
# across multiple sheets, reflecting a fresh scheduled-runner data pull.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# @@ -2540,22 +2540,22 @@
$ws.Range("H39").Value = 11111196
$ws.Range("I39").Value = 12500083
$ws.Range("K39").Value = 37500249
$ws.Range("M39").Value = -37499953

# @@ -3694,25 +3694,25 @@
$ws.Range("H62").Value = 9855.714
$ws.Range("I62").Value = 8000
$ws.Range("J62").Value = 12330
$ws.Range("K62").Value = 8000
$ws.Range("L62").Value = 12330
$ws.Range("M62").Value = -7376
$ws.Range("N62").Value = -13578

# @@ -3847,25 +3847,25 @@
$ws.Range("H65").Value = 9855.714
$ws.Range("I65").Value = 8000
$ws.Range("J65").Value = 12330
$ws.Range("K65").Value = 40000
$ws.Range("L65").Value = 61650
$ws.Range("M65").Value = -36880
$ws.Range("N65").Value = -67890

# @@ -6262,19 +6262,22 @@
$ws.Range("H113").Value = 1600
$ws.Range("J113").Value = 1600
$ws.Range("L113").Value = 1600
$ws.Range("N113").Value = -8108

# @@ -7658,22 +7661,22 @@
$ws.Range("H141").Value = 5421.6924
$ws.Range("I141").Value = 5317
$ws.Range("K141").Value = 15951
$ws.Range("M141").Value = -10771


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# @@ -9271,25 +9274,25 @@
$ws.Range("H32").Value = 5170.353
$ws.Range("I32").Value = 4931
$ws.Range("J32").Value = 9000
$ws.Range("K32").Value = 4931
$ws.Range("L32").Value = 9000
$ws.Range("M32").Value = -4644
$ws.Range("N32").Value = -9574

# @@ -11302,22 +11305,22 @@
$ws.Range("H74").Value = 1865.2
$ws.Range("I74").Value = 1922.0834
$ws.Range("K74").Value = 1922.0834
$ws.Range("M74").Value = -1048.0834

# @@ -11449,22 +11452,22 @@
$ws.Range("H77").Value = 1865.2
$ws.Range("I77").Value = 1922.0834
$ws.Range("K77").Value = 9610.416999999999
$ws.Range("M77").Value = -5242.416999999999

# @@ -12420,22 +12423,22 @@
$ws.Range("H97").Value = 1918.6923
$ws.Range("I97").Value = 1585.7273
$ws.Range("K97").Value = 1585.7273
$ws.Range("M97").Value = -1089.7273

# @@ -14153,22 +14156,22 @@
$ws.Range("H132").Value = 280389.25
$ws.Range("I132").Value = 324920.6
$ws.Range("K132").Value = 974761.7999999999
$ws.Range("M132").Value = -972231.7999999999


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# @@ -15607,25 +15610,25 @@
$ws.Range("H20").Value = 3835.6365
$ws.Range("I20").Value = 3437.6
$ws.Range("J20").Value = 4167.3335
$ws.Range("K20").Value = 3437.6
$ws.Range("L20").Value = 4167.3335
$ws.Range("M20").Value = -3190.6
$ws.Range("N20").Value = -4661.3335

# @@ -18817,23 +18820,20 @@
$ws.Range("H86").Value = 911.75
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

# @@ -18967,23 +18967,20 @@
$ws.Range("H89").Value = 911.75
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# @@ -19766,22 +19763,22 @@
$ws.Range("H105").Value = 3215.1428
$ws.Range("I105").Value = 1733.6
$ws.Range("K105").Value = 1733.6
$ws.Range("M105").Value = 13.40000000000009

# @@ -19867,22 +19864,22 @@
$ws.Range("H107").Value = 2470.7334
$ws.Range("I107").Value = 1451.2222
$ws.Range("K107").Value = 1451.2222
$ws.Range("M107").Value = 468.7778000000001

# @@ -20801,22 +20798,22 @@
$ws.Range("H126").Value = 116999
$ws.Range("J126").Value = 116999
$ws.Range("L126").Value = 116999
$ws.Range("N126").Value = -126879

# @@ -21098,22 +21095,22 @@
$ws.Range("H132").Value = 126992.6
$ws.Range("J132").Value = 126992.6
$ws.Range("L132").Value = 126992.6
$ws.Range("N132").Value = -137112.6

# @@ -21196,22 +21193,22 @@
$ws.Range("H134").Value = 11496656
$ws.Range("I134").Value = 1943.2174
$ws.Range("K134").Value = 5829.6522
$ws.Range("M134").Value = -3294.6522


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# @@ -22080,25 +22077,25 @@
$ws.Range("H10").Value = 1519.8948
$ws.Range("I10").Value = 1703.9166
$ws.Range("J10").Value = 1204.4286
$ws.Range("K10").Value = 1703.9166
$ws.Range("L10").Value = 1204.4286
$ws.Range("M10").Value = -1564.9166
$ws.Range("N10").Value = -1482.4286

# @@ -23118,25 +23115,25 @@
$ws.Range("H31").Value = 2947.5854
$ws.Range("I31").Value = 2094.577
$ws.Range("J31").Value = 4426.1333
$ws.Range("K31").Value = 2094.577
$ws.Range("L31").Value = 4426.1333
$ws.Range("M31").Value = -1799.577
$ws.Range("N31").Value = -5016.1333

# @@ -23271,25 +23268,25 @@
$ws.Range("H34").Value = 2947.5854
$ws.Range("I34").Value = 2094.577
$ws.Range("J34").Value = 4426.1333
$ws.Range("K34").Value = 2094.577
$ws.Range("L34").Value = 4426.1333
$ws.Range("M34").Value = -1892.577
$ws.Range("N34").Value = -4830.1333

# @@ -28097,22 +28094,22 @@
$ws.Range("H132").Value = 3399.2144
$ws.Range("I132").Value = 2781.0908
$ws.Range("K132").Value = 8343.2724
$ws.Range("M132").Value = -5813.2724

# @@ -28195,22 +28192,22 @@
$ws.Range("H134").Value = 2118.0476
$ws.Range("I134").Value = 1521.9412
$ws.Range("K134").Value = 4565.8236
$ws.Range("M134").Value = -2030.8236


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# @@ -41756,22 +41753,22 @@
$ws.Range("H126").Value = 4408.75
$ws.Range("I126").Value = 3829.5386
$ws.Range("K126").Value = 11488.6158
$ws.Range("M126").Value = -9018.6158

# @@ -42053,22 +42050,22 @@
$ws.Range("H132").Value = 3080.7334
$ws.Range("I132").Value = 3080.7334
$ws.Range("K132").Value = 9242.200199999999
$ws.Range("M132").Value = -6712.200199999999


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# @@ -43320,25 +43317,25 @@
$ws.Range("H16").Value = 1738.3889
$ws.Range("I16").Value = 1799.4
$ws.Range("J16").Value = 1433.3334
$ws.Range("K16").Value = 1799.4
$ws.Range("L16").Value = 1433.3334
$ws.Range("M16").Value = -1629.4
$ws.Range("N16").Value = -1773.3334

# @@ -44407,23 +44404,20 @@
$ws.Range("H38").Value = 30000
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

# @@ -49001,22 +48995,22 @@
$ws.Range("H132").Value = 359916.03
$ws.Range("I132").Value = 457329.78
$ws.Range("K132").Value = 1371989.34
$ws.Range("M132").Value = -1369459.34


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# @@ -49931,25 +49925,22 @@
$ws.Range("H9").Value = 2412.5
$ws.Range("I9").Value = 2412.5
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 2412.5
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -2272.5
$ws.Range("N9").ClearContents()

# @@ -52001,25 +51992,19 @@
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

# @@ -52053,25 +52038,22 @@
$ws.Range("H52").Value = 59993
$ws.Range("I52").Value = 59993
$ws.Range("J52").Value = 59993
$ws.Range("K52").Value = 59993
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -59767
$ws.Range("N52").ClearContents()

# @@ -54405,25 +54387,25 @@
$ws.Range("H100").Value = 1429.6666
$ws.Range("I100").Value = 1708.5
$ws.Range("J100").Value = 872
$ws.Range("K100").Value = 3417
$ws.Range("L100").Value = 1744
$ws.Range("M100").Value = -2876
$ws.Range("N100").Value = -2826

# @@ -55679,22 +55661,22 @@
$ws.Range("I126").Value = 2987.5
$ws.Range("K126").Value = 8962.5


Write-Host "Applied scheduled-runner market data update."